# The "Requisitos" section at the end of the document used to be followed
# by a blank paragraph, a "Ver no Jupiter ..." line and a "(c) 2020 ..."
# footer line. The edit removes those three trailing paragraphs, leaving
# the "LOM3005: ..." line followed directly by the blank paragraph / page
# break that already existed further down.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph via Find (robust to absolute
# paragraph-index assumptions).
$hit = $d.Content
$found = $hit.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the 'Ver no Jupiter' paragraph"
}

# Resolve the matching Paragraph object in the document's Paragraphs
# collection so we can walk to its neighbours.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $hit.Start) {
        $target = $p
        break
    }
}

if ($null -eq $target) {
    throw "Could not resolve paragraph object for the matched text"
}

# The paragraph immediately before it is the empty spacer paragraph, and
# the one immediately after is the "(c) 2020 ..." footer line - both are
# removed along with the "Ver no Jupiter ..." paragraph itself.
$prevPara = $target.Previous()
$nextPara = $target.Next()

$deleteStart = $prevPara.Range.Start
$deleteEnd = $nextPara.Range.End

$d.Range($deleteStart, $deleteEnd).Delete()
